# Updated cryptos list on Mon May  6 16:17:35 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that must stay TEXT even when they look like
# numbers (e.g. "589.11" or thousand-grouped "63.712.74"). Force text format
# before writing so Excel does not auto-convert them to numeric values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.712.74"
$ws.Range("E2").Value = "  -1.34%  "

$ws.Range("D3").Value = "3.089.32"
$ws.Range("E3").Value = "  -2.36%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "589.11"
$ws.Range("E5").Value = "  -0.92%  "

$ws.Range("D6").Value = "153.46"
$ws.Range("E6").Value = "  +4.16%  "

$ws.Range("D8").Value = "0.548"
$ws.Range("E8").Value = "  +2.82%  "

$ws.Range("D9").Value = "3.084.67"
$ws.Range("E9").Value = "  -2.47%  "

$ws.Range("E10").Value = "  -2.36%  "

$ws.Range("D11").Value = "5.89"
$ws.Range("E11").Value = "  -0.73%  "

$ws.Range("D12").Value = "0.461"
$ws.Range("E12").Value = "  -1.07%  "

$ws.Range("D13").Value = "37.99"
$ws.Range("E13").Value = "  +1.23%  "

$ws.Range("D14").Value = "0.0000243"
$ws.Range("E14").Value = "  -2.63%  "

$ws.Range("D15").Value = "3.597.27"
$ws.Range("E15").Value = "  -2.31%  "

$ws.Range("E16").Value = "  -1.95%  "

$ws.Range("E17").Value = "  -1.25%  "

$ws.Range("D18").Value = "63.675.32"
$ws.Range("E18").Value = "  -0.88%  "

$ws.Range("D19").Value = "3.085.10"
$ws.Range("E19").Value = "  -2.42%  "

$ws.Range("D20").Value = "473.50"
$ws.Range("E20").Value = "  +0.26%  "

$ws.Range("D21").Value = "14.75"
$ws.Range("E21").Value = "  +0.82%  "

$ws.Range("D22").Value = "0.724"
$ws.Range("E22").Value = "  -2.20%  "

$ws.Range("E23").Value = "  -0.43%  "

$ws.Range("E24").Value = "  +3.73%  "

$ws.Range("D25").Value = "13.19"
$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("D26").Value = "81.87"
$ws.Range("E26").Value = "  +0.17%  "

$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("D28").Value = "9.89"
$ws.Range("E28").Value = "  +2.28%  "

$ws.Range("E29").Value = "  -1.75%  "

$ws.Range("E30").Value = "  -2.34%  "

$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.14%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "2.21"
$ws.Range("E32").Value = "  -1.95%  "

$ws.Range("E33").Value = "  +4.15%  "

$ws.Range("D34").Value = "27.44"
$ws.Range("E34").Value = "  -0.72%  "

$ws.Range("D35").Value = "0.0₃0853"
$ws.Range("E35").Value = "  -0.30%  "

$ws.Range("D36").Value = "1.06"
$ws.Range("E36").Value = "  -1.24%  "

$ws.Range("D37").Value = "3.42"
$ws.Range("E37").Value = "  +4.53%  "

$ws.Range("E38").Value = "  -2.22%  "

$ws.Range("D39").Value = "2.25"
$ws.Range("E39").Value = "  -4.25%  "

$ws.Range("D40").Value = "9.37"
$ws.Range("E40").Value = "  +1.63%  "

$ws.Range("D41").Value = "50.71"
$ws.Range("E41").Value = "  -2.36%  "

$ws.Range("D42").Value = "451.51"
$ws.Range("E42").Value = "  -1.30%  "

$ws.Range("D43").Value = "0.287"
$ws.Range("E43").Value = "  -2.61%  "

$ws.Range("E44").Value = "  -2.32%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.837.42"
$ws.Range("E45").Value = "  -3.57%  "

$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").Value = "40.02"
$ws.Range("E46").Value = "  -1.78%  "

$ws.Range("E47").Value = "  -0.79%  "

$ws.Range("D48").Value = "130.46"
$ws.Range("E48").Value = "  +0.86%  "

$ws.Range("D49").Value = "25.73"
$ws.Range("E49").Value = "  +4.66%  "

$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.06%  "

$ws.Range("D51").Value = "2.28"
$ws.Range("E51").Value = "  +0.52%  "
